# Replace COSOPT's MeanExpLev-derived "Expressed" counts with counts
# derived from the raw expression-count data (commit: "Determine whether
# a gene is expressed using expression counts data, not COSOPT's
# MeanExpLev").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> (Rhythmic, Expressed)
$updates = @{
    2 = @(22333, 35205)
    3 = @(22383, 35299)
    4 = @(25661, 35491)
    5 = @(19100, 35013)
    6 = @(24568, 35239)
    7 = @(26490, 35492)
    8 = @(19064, 35005)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
}

# The author's last selection before saving, per the sheetView diff.
$ws.Range("F20").Select()
